$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text values (shared strings) that visibly changed
$ws.Range("B10").Value = "7692.93 ± 87.6281"
$ws.Range("B14").Value = "15316.7 ± 5605.41"

# Update numeric values
$ws.Range("B12").Value = 7443
$ws.Range("B13").Value = 0.2
$ws.Range("B15").Value = 9100
$ws.Range("B16").Value = 22000

# Update the active selection/cell
$ws.Range("D19").Select()
